$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the swap between rows, by letter -> index
# A=1 B=2 E=5 F=6 G=7 H=8 M=13 Q=17 R=18

function Swap-Cells($row1, $row2, $col) {
    $c1 = $ws.Cells.Item($row1, $col)
    $c2 = $ws.Cells.Item($row2, $col)
    $v1 = $c1.Value()
    $v2 = $c2.Value()
    $c1.Value = $v2
    $c2.Value = $v1
}

# Swap rows 6 and 8 for columns A, B, E, F, G, H, M, Q, R
$cols1 = 1, 2, 5, 6, 7, 8, 13, 17, 18
foreach ($col in $cols1) {
    Swap-Cells 6 8 $col
}

# Swap rows 10 and 11 for columns A, Q, R
$cols2 = 1, 17, 18
foreach ($col in $cols2) {
    Swap-Cells 10 11 $col
}
